# TreeSeqCompTime.xlsx update
# - Adds two new parameter rows (A=15, A=20) ahead of the existing series
#   and three more (A=30,35,40) in place of the old empty placeholder rows
#   (A=750,1000,2000,3000), i.e. updates the parameter space.
# - Right-aligns the new "Peak Memory_GB" (D) values for the middle block
#   of rows with a new cell style (Menlo font, right aligned).
# - Renames the sheet to match the exported CSV name and tweaks the saved
#   window geometry / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- make room: insert one row so the data block grows from 10 to 12 rows ---
$ws.Range("A2:F2").Insert(-4121)

# --- helper style constants ---
$xlRight = -4152

function Set-RowStyle($row, $dStyleKind) {
    # Column A: plain data style (Menlo Regular 11, default color)
    $a = $ws.Cells.Item($row, 1)
    $a.Font.Name = "Menlo Regular"
    $a.Font.Size = 11

    # Column E/F keep their existing look, but make sure E matches A's style
    $e = $ws.Cells.Item($row, 5)
    $e.Font.Name = "Menlo Regular"
    $e.Font.Size = 11

    $f = $ws.Cells.Item($row, 6)
    $f.Font.Name = "Menlo Regular"
    $f.Font.Size = 12
}

# Row data: A (CPU), B (Hours raw seconds), D (Peak Memory_GB)
$rows = @(
    @{ R=2;  A=15;  B=33315.4; D=7.3789999999999996;  BStyle="menlo-black";   DStyle="menloreg-black-right" },
    @{ R=3;  A=20;  B=30716.6; D=9.3859999999999992;  BStyle="menlo-black";   DStyle="menloreg-black-right" },
    @{ R=4;  A=25;  B=26024.2; D=10.143000000000001;  BStyle="menlo-black";   DStyle="menlo-black-right" },
    @{ R=5;  A=30;  B=26777.8; D=12.818;               BStyle="menlo-black";   DStyle="menlo-black-right" },
    @{ R=6;  A=35;  B=24486;   D=13.61;                BStyle="menlo-black";   DStyle="menlo-black-right" },
    @{ R=7;  A=40;  B=23265.3; D=15.286;                BStyle="menlo-black";   DStyle="menlo-black-right" },
    @{ R=8;  A=50;  B=20554.7; D=18.760000000000002;   BStyle="menlo-black";   DStyle="menlo-black-right" },
    @{ R=9;  A=100; B=19575.7; D=35.820999999999998;   BStyle="menloreg-black";DStyle="menloreg-black-right" },
    @{ R=10; A=250; B=21464.5; D=87.296999999999997;   BStyle="menloreg-black";DStyle="menloreg-black-right" },
    @{ R=11; A=500; B=21371;   D=173.13399999999999;   BStyle="menloreg-theme";DStyle="menloreg-black-right" }
)

foreach ($row in $rows) {
    $r = $row.R

    # Column A
    $cA = $ws.Cells.Item($r, 1)
    $cA.Value2 = $row.A
    $cA.Font.Name = "Menlo Regular"
    $cA.Font.Size = 11

    # Column B
    $cB = $ws.Cells.Item($r, 2)
    $cB.Value2 = $row.B
    if ($row.BStyle -eq "menlo-black") {
        $cB.Font.Name = "Menlo"
        $cB.Font.Size = 11
        $cB.Font.Color = 0
    } elseif ($row.BStyle -eq "menloreg-black") {
        $cB.Font.Name = "Menlo Regular"
        $cB.Font.Size = 11
        $cB.Font.Color = 0
    } else {
        $cB.Font.Name = "Menlo Regular"
        $cB.Font.Size = 11
    }

    # Column C: formula, value computed by recalculation
    $cC = $ws.Cells.Item($r, 3)
    $cC.Formula = "=B" + $r + "/3600"
    $cC.Font.Name = "Menlo Regular"
    $cC.Font.Size = 11
    $cC.NumberFormat = "0.00"

    # Column D
    $cD = $ws.Cells.Item($r, 4)
    $cD.Value2 = $row.D
    if ($row.DStyle -eq "menlo-black-right") {
        $cD.Font.Name = "Menlo"
        $cD.Font.Size = 11
        $cD.Font.Color = 0
    } else {
        $cD.Font.Name = "Menlo Regular"
        $cD.Font.Size = 11
    }
    $cD.HorizontalAlignment = $xlRight

    # Column E / F
    $cE = $ws.Cells.Item($r, 5)
    $cE.Font.Name = "Menlo Regular"
    $cE.Font.Size = 11

    $cF = $ws.Cells.Item($r, 6)
    $cF.Font.Name = "Menlo Regular"
    $cF.Font.Size = 12
}

# Shared formula for the contiguous fill-down block C3:C7 (rows for A=20..40)
$ws.Range("C3:C7").Formula = "=B3/3600"
foreach ($r in 3..7) {
    $ws.Cells.Item($r, 3).Font.Name = "Menlo Regular"
    $ws.Cells.Item($r, 3).Font.Size = 11
    $ws.Cells.Item($r, 3).NumberFormat = "0.00"
}

# --- selection / view state ---
$ws.Range("B9").Select()

# --- sheet name matches the exported CSV source ---
$ws.Name = "TreeSeqCompTime.csv"

# --- saved window geometry ---
$excel.ActiveWindow.WindowState = -4143
$excel.Width = 14300
$excel.Left = 13240
